$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet updates (Version, Status, Date, Contact)
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.4.0-snapshot-1"
$meta.Range("B6").Value = "draft"
$meta.Range("B8").Value = "2024-05-23T12:16:26+00:00"
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"

# ---------------------------------------------------------------------------
# 2. Elements sheet - swap the "Mapping" columns AK (37) and AL (38):
#    header text, per-row data, and column widths all need to be swapped.
# ---------------------------------------------------------------------------
$els = $wb.Worksheets.Item("Elements")

# Swap header labels in row 1
$akHeader = $els.Cells.Item(1, 37).Value()
$alHeader = $els.Cells.Item(1, 38).Value()
$els.Cells.Item(1, 37).Value = $alHeader
$els.Cells.Item(1, 38).Value = $akHeader

# Swap the data cells for each data row (rows 2-6)
$lastRow = 6
for ($r = 2; $r -le $lastRow; $r++) {
    $akCell = $els.Cells.Item($r, 37)
    $alCell = $els.Cells.Item($r, 38)
    $akVal = $akCell.Value()
    $alVal = $alCell.Value()
    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Swap the column widths so the wider column lines up with the longer text
$els.Columns.Item(37).ColumnWidth = 64.0
$els.Columns.Item(38).ColumnWidth = 24.166666666666668
